# This script "re-distributes" the per-record data found in rows 11-64 of the
# Artfynd worksheet to a different set of rows (a row permutation keyed by the
# unique record Id in column A), matching a refreshed export of the same
# underlying observations. While doing so it also:
#   - rounds the Ost/Nord coordinate columns (Q/R) to whole numbers
#   - drops the Starttid/Sluttid (Z/AB) time-of-day values (always "00:00")
#
# Because every row's data moves to a new row, we first snapshot the full
# current contents of rows 11-64 (columns A..AY) into memory, then write the
# snapshot back out according to the row-to-row mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 11
$lastRow = 64
$firstCol = 1    # A
$lastCol = 51    # AY

$colQ = 17   # Ost
$colR = 18   # Nord
$colZ = 26   # Starttid
$colAB = 28  # Sluttid

# Target row -> source row (both refer to *original* row positions before
# this script runs). This is a fixed permutation of rows 11..64 derived from
# matching each record's unique Id (column A) between the before/after data.
$rowSourceMap = @{}
$rowSourceMap[11] = 18
$rowSourceMap[12] = 23
$rowSourceMap[13] = 20
$rowSourceMap[14] = 53
$rowSourceMap[15] = 55
$rowSourceMap[16] = 44
$rowSourceMap[17] = 39
$rowSourceMap[18] = 13
$rowSourceMap[19] = 46
$rowSourceMap[20] = 35
$rowSourceMap[21] = 40
$rowSourceMap[22] = 62
$rowSourceMap[23] = 51
$rowSourceMap[24] = 48
$rowSourceMap[25] = 57
$rowSourceMap[26] = 15
$rowSourceMap[27] = 12
$rowSourceMap[28] = 59
$rowSourceMap[29] = 32
$rowSourceMap[30] = 34
$rowSourceMap[31] = 60
$rowSourceMap[32] = 26
$rowSourceMap[33] = 21
$rowSourceMap[34] = 33
$rowSourceMap[35] = 52
$rowSourceMap[36] = 14
$rowSourceMap[37] = 47
$rowSourceMap[38] = 22
$rowSourceMap[39] = 27
$rowSourceMap[40] = 29
$rowSourceMap[41] = 16
$rowSourceMap[42] = 38
$rowSourceMap[43] = 45
$rowSourceMap[44] = 43
$rowSourceMap[45] = 28
$rowSourceMap[46] = 61
$rowSourceMap[47] = 63
$rowSourceMap[48] = 42
$rowSourceMap[49] = 56
$rowSourceMap[50] = 25
$rowSourceMap[51] = 37
$rowSourceMap[52] = 64
$rowSourceMap[53] = 49
$rowSourceMap[54] = 17
$rowSourceMap[55] = 41
$rowSourceMap[56] = 58
$rowSourceMap[57] = 36
$rowSourceMap[58] = 50
$rowSourceMap[59] = 31
$rowSourceMap[60] = 19
$rowSourceMap[61] = 11
$rowSourceMap[62] = 54
$rowSourceMap[63] = 24
$rowSourceMap[64] = 30

# 1) Snapshot every cell value of rows 11..64 before making any changes.
$cache = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2)
    }
    $cache[$r] = $rowVals
}

# 2) Write the snapshot back out row by row according to the mapping,
#    applying the Q/R rounding and clearing Z/AB along the way.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowSourceMap[$r]
    $srcVals = $cache[$srcRow]

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $val = $srcVals[$c - $firstCol]
        $cell = $ws.Cells.Item($r, $c)

        if ($c -eq $colZ -or $c -eq $colAB) {
            # Starttid/Sluttid are dropped entirely in the refreshed export.
            $cell.ClearContents() | Out-Null
            continue
        }

        if ($c -eq $colQ -or $c -eq $colR) {
            if ($val -eq $null) {
                $cell.ClearContents() | Out-Null
            } else {
                $cell.Value = [Math]::Round([double]$val)
            }
            continue
        }

        if ($val -eq $null) {
            $cell.ClearContents() | Out-Null
        } else {
            if ($val -is [string]) {
                # Force plain text so values such as "2023-09-05" are kept as
                # literal strings instead of being auto-converted to dates.
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }
}

Write-Host "Row permutation applied to rows $firstRow-$lastRow."
